$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 / 18: WrappedBTC and WrappedEther swap places (name + link), plus new price/volume values ---
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"

# --- Price (D) and Volume(1h) (E) updates across all rows ---
# Cells whose new text looks like a plain number (e.g. "603.41") need to be forced
# to remain text (matching the source data, which stores these as inline strings),
# otherwise Excel auto-converts them to numeric values. We briefly mark the cell as
# Text format, assign the value, then clear formatting again so the cell keeps the
# workbook default style (no explicit style index), just like the other text cells.
$ws.Range("D2").Value = "70.457.19"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.624.24"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.12%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "603.41"
$r.ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "196.68"
$r.ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.91%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.647"
$r.ClearFormats()
$ws.Range("E10").Value = "  -0.96%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "53.38"
$r.ClearFormats()
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  +0.02%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "9.56"
$r.ClearFormats()
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "4.195.82"
$ws.Range("E14").Value = "  +2.33%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "598.63"
$r.ClearFormats()
$ws.Range("E15").Value = "  -1.16%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "13.02"
$r.ClearFormats()
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "3.645.76"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "70.533.75"
$ws.Range("E18").Value = "  +0.57%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "19.03"
$r.ClearFormats()
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E21").Value = "  -0.21%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "18.07"
$r.ClearFormats()
$ws.Range("E22").Value = "  -0.04%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "5.19"
$r.ClearFormats()
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  -0.61%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "4.62"
$r.ClearFormats()
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -3.53%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "10.72"
$r.ClearFormats()
$ws.Range("E27").Value = "  -2.07%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "9.68"
$r.ClearFormats()
$ws.Range("E28").Value = "  +0.50%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "33.82"
$r.ClearFormats()
$ws.Range("E29").Value = "  +0.53%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "4.62"
$r.ClearFormats()
$ws.Range("E30").Value = "  +6.46%  "
$ws.Range("E31").Value = "  +1.11%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "12.27"
$r.ClearFormats()
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("E33").Value = "  +3.17%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "63.49"
$r.ClearFormats()
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "0.0₃0888"
$ws.Range("E35").Value = "  +4.16%  "
$ws.Range("D36").Value = "3.913.29"
$ws.Range("E36").Value = "  +3.73%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "538.53"
$r.ClearFormats()
$ws.Range("E37").Value = "  +9.79%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.09"
$r.ClearFormats()
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  -0.03%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "36.94"
$r.ClearFormats()
$ws.Range("E40").Value = "  +0.59%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.389"
$r.ClearFormats()
$ws.Range("E41").Value = "  -1.44%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "3.54"
$r.ClearFormats()
$ws.Range("E42").Value = "  -3.45%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.0465"
$r.ClearFormats()
$ws.Range("E43").Value = "  +1.57%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.134"
$r.ClearFormats()
$ws.Range("E44").Value = "  -0.02%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.42"
$r.ClearFormats()
$ws.Range("E45").Value = "  +3.47%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "2.88"
$r.ClearFormats()
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  -0.18%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.000248"
$r.ClearFormats()
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  +1.43%  "
